# Add new power plant / electricity source rows to the BGDPbES sheet
# (issues #280 and #99): hard coal w CCS, natural gas combined cycle w CCS,
# biomass w CCS, lignite w CCS, small modular reactor, hydrogen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

$newSources = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
for ($i = 0; $i -lt $newSources.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newSources[$i]
    for ($col = 2; $col -le 37; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}

# Update sheet selection/activation state: BGDPbES is no longer the tab that
# is active/selected; "About" becomes the active sheet, while BGDPbES keeps
# a remembered selection at the first empty row below the new data (A25).
$ws.Range("A25").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
